# Minor changes to figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two bottom-most data rows (9 and 10) so the table shrinks from
# 9 data rows (A1:H10) down to 7 data rows (A1:H8).
$ws.Rows(9).Delete()
$ws.Rows(9).Delete()

# The refreshed figures no longer carry the DiffDE/AgeCompDE/RateCompDE/
# relAgeDE/relRateDE breakdown columns - clear D:H for every remaining
# data row.
$ws.Range("D2:H8").ClearContents()

# Row 2 - SouthKorea
$ws.Range("A2").Value = "SouthKorea"
$ws.Range("B2").Value = 43942
$ws.Range("C2").Value = 0.02218477955630441

# Row 3 - China
$ws.Range("A3").Value = "China"
$ws.Range("B3").Value = 43872
$ws.Range("C3").Value = 0.02290025071633238

# Row 4 - France
$ws.Range("A4").Value = "France"
$ws.Range("B4").Value = 43914
$ws.Range("C4").Value = 0.03983587515221891

# Row 5 - USA
$ws.Range("A5").Value = "USA"
$ws.Range("B5").Value = 43940
$ws.Range("C5").Value = 0.06870385174884934

# Row 6 - Spain
$ws.Range("A6").Value = "Spain"
$ws.Range("B6").Value = 43937
$ws.Range("C6").Value = 0.1050210003716739

# Row 7 - Italy
$ws.Range("A7").Value = "Italy"
$ws.Range("B7").Value = 43941
$ws.Range("C7").Value = 0.1272752828730058

# Row 8 - Germany (figure for CFR2 not yet available this refresh)
$ws.Range("A8").Value = "Germany"
$ws.Range("B8").Value = 43941
$ws.Range("C8").ClearContents()
